$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Rewrite the bullets under "KEY ACHIEVEMENTS AND IMPACT" so they read as
# punchy, impact-focused accomplishment statements instead of job-duty
# descriptions. Scope everything to that section so the near-duplicate
# bullets living under "PROFESSIONAL EXPERIENCE" (Siege Analytics) are left
# untouched.
# ---------------------------------------------------------------------------

function Find-ParaIndex($startIdx, $endIdx, $pattern) {
    for ($i = $startIdx; $i -le $endIdx; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

$count = $d.Paragraphs.Count

$sectionStart = Find-ParaIndex 1 $count "KEY ACHIEVEMENTS AND IMPACT"
if ($sectionStart -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT section"
}

$sectionEnd = Find-ParaIndex ($sectionStart + 1) $count "TECHNICAL SKILLS"
if ($sectionEnd -eq -1) {
    $sectionEnd = $count
}

# The "FEC analysis" bullet is dropped entirely in the new version -- delete
# its whole paragraph (including the paragraph mark) first so the other
# bullets' indices above it stay valid while we edit them next.
$fecIdx = Find-ParaIndex $sectionStart $sectionEnd "Built real-time FEC analysis systems"
if ($fecIdx -eq -1) {
    throw "Could not find the FEC analysis bullet to remove"
}
$d.Paragraphs.Item($fecIdx).Range.Delete()

# Recompute the section bounds post-deletion and rewrite the remaining
# bullets in place, preserving their paragraph/run formatting.
$count = $d.Paragraphs.Count
$sectionEnd = Find-ParaIndex $sectionStart $count "TECHNICAL SKILLS"
if ($sectionEnd -eq -1) {
    $sectionEnd = $count
}

$replacements = @(
    @{ Find = "Discovered systematic race coding errors affecting all Black and Asian-American voters"; Replace = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions" },
    @{ Find = "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%"; Replace = "• 178% accuracy improvement in racial classification algorithms" },
    @{ Find = "Built redistricting platform used by thousands of analysts nationwide"; Replace = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%" },
    @{ Find = "Achieved 87% prediction accuracy for voter turnout"; Replace = "• `$4.7M savings enabled nonprofit access" },
    @{ Find = "Provided expert testimony and press briefings on electoral data integrity"; Replace = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations" }
)

foreach ($r in $replacements) {
    $idx = Find-ParaIndex $sectionStart $sectionEnd $r.Find
    if ($idx -eq -1) {
        throw "Could not find bullet matching '$($r.Find)'"
    }
    $d.Paragraphs.Item($idx).Range.Text = $r.Replace
}
